$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing weekly price row (row 59, date 2020-12-24) needs to be kept as
# historical data, but shifted down to make room for a new weekly entry that
# reuses row 59 with fresh data (date 2022-01-24).
#
# So: insert a new blank row at 60 (this pushes the old row 60 down to 61,
# while row 59 is untouched), copy the current (soon to be historical)
# contents of row 59 into the newly inserted row 60, and then overwrite row
# 59 with the new week's figures.

$ws.Rows("60:60").Insert()

$ws.Range("A59:R59").Copy() | Out-Null
$ws.Range("A60").PasteSpecial(-4104) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("D59").Value2 = 44585
$ws.Range("J59").Value2 = 25
$ws.Range("K59").Value2 = 28000
$ws.Range("L59").Value2 = 28000
$ws.Range("M59").Value2 = 28000
$ws.Range("P59").Value2 = 1120
